$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 7 - "VR Home Design" (new Title Slide, same layout as slide 1)
# ---------------------------------------------------------------------------
$s7 = $p.Slides.Add(7, 1)
$s7.Shapes.Item(1).TextFrame.TextRange.Text = "VR Home Design"

# ---------------------------------------------------------------------------
# Slide 8 - "Description"
# ---------------------------------------------------------------------------
$s8 = $p.Slides.Add(8, 2)
$s8.Shapes.Item(1).TextFrame.TextRange.Text = "Description"
$tf8 = $s8.Shapes.Item(2).TextFrame
$tf8.TextRange.Text = "A VR application that can assist in the basic design layout and décor of a home.`r`rPlace décor and roam the rooms of a house of your making."
$tf8.TextRange.Paragraphs(2, 1).ParagraphFormat.Bullet.Visible = 0

# ---------------------------------------------------------------------------
# Slide 9 - "Core Features"
# ---------------------------------------------------------------------------
$s9 = $p.Slides.Add(9, 2)
$s9.Shapes.Item(1).TextFrame.TextRange.Text = "Core Features"
$tf9 = $s9.Shapes.Item(2).TextFrame
$tf9.TextRange.Text = "Move around the home in VR`r`rPlace and resize rooms`r`rSave your work`r`rAdd basic furnishing and décor to the room.`r`rDesign within VR"
$tf9.TextRange.Paragraphs(2, 1).ParagraphFormat.Bullet.Visible = 0
$tf9.TextRange.Paragraphs(4, 1).ParagraphFormat.Bullet.Visible = 0
$tf9.TextRange.Paragraphs(6, 1).ParagraphFormat.Bullet.Visible = 0

# ---------------------------------------------------------------------------
# Slide 10 - "Stretch Goals"
# ---------------------------------------------------------------------------
$s10 = $p.Slides.Add(10, 2)
$s10.Shapes.Item(1).TextFrame.TextRange.Text = "Stretch Goals"
$tf10 = $s10.Shapes.Item(2).TextFrame
$tf10.TextRange.Text = "Add custom object models as décor and furnishing`r`rConnect multiple users to roam among the same simulated home.`r`rDesign templates`r`rUse of various building materials`r`rAdvanced lighting simulation"

# ---------------------------------------------------------------------------
# Slide 11 - "Pros"
# ---------------------------------------------------------------------------
$s11 = $p.Slides.Add(11, 2)
$s11.Shapes.Item(1).TextFrame.TextRange.Text = "Pros"
$tf11 = $s11.Shapes.Item(2).TextFrame
$tf11.TextRange.Text = "Game engines with built in VR support can be used for VR rendering assistance.`r`rMobile VR is accessible`r`rBasic models of necessary items are obtainable"

# ---------------------------------------------------------------------------
# Slide 12 - "Cons"
# ---------------------------------------------------------------------------
$s12 = $p.Slides.Add(12, 2)
$s12.Shapes.Item(1).TextFrame.TextRange.Text = "Cons"
$tf12 = $s12.Shapes.Item(2).TextFrame
$tf12.TextRange.Text = "Computationally intensive`r`rBest method for input?`r`rNeeds many accessible models"
